# Apply the cryptos-list price/volume refresh captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.586.04"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "2.164.38"
$ws.Range("E3").Value = "  -3.25%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.26"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.91"
$ws.Range("E7").Value = "  -3.89%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.576"
$ws.Range("E9").Value = "  -5.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.77"
$ws.Range("E10").Value = "  -7.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0904"
$ws.Range("E11").Value = "  -5.99%  "
$ws.Range("E12").Value = "  -5.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.100"
$ws.Range("E13").Value = "  -3.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.69"
$ws.Range("E14").Value = "  -4.69%  "
$ws.Range("D15").Value = "2.489.49"
$ws.Range("E15").Value = "  -3.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.13"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").Value = "2.152.62"
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.780"
$ws.Range("E18").Value = "  -7.34%  "
$ws.Range("D19").Value = "41.448.86"
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.82"
$ws.Range("E21").Value = "  -4.50%  "
$ws.Range("E22").Value = "  -7.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.82"
$ws.Range("E23").Value = "  -12.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "226.84"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  -7.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.01"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.76"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.81"
$ws.Range("E33").Value = "  +9.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0769"
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.13"
$ws.Range("E35").Value = "  -9.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.121"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  -5.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0303"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.10"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.95"
$ws.Range("E41").Value = "  -11.44%  "
$ws.Range("E42").Value = "  -6.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.14"
$ws.Range("E43").Value = "  -9.62%  "
$ws.Range("E44").Value = "  -5.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.41"
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("E46").Value = "  -4.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.03"
$ws.Range("E47").Value = "  -7.39%  "
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("E49").Value = "  -5.42%  "
$ws.Range("E50").Value = "  -8.01%  "
$ws.Range("E51").Value = "  -2.46%  "
